$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Salutation header block ("{{identity.localized_salutation}} {{identity.
#    localized_title}}") -> single merge field "{{identity.address_block}}"
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2).Range
$p2.Find.Execute("{{identity.localized_salutation}} {{identity.localized_title}}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{{identity.address_block}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Drop the whole conditional name / PO-box / street / country address
#    block (old paragraphs 3-11): first empty the first of those paragraphs
#    (it stays behind as a blank line) and then delete paragraphs 4-12
#    outright (content + paragraph marks) so they collapse away entirely.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute("{{identity.first_name}} {{identity.last_name}}", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$rangeStart = $d.Paragraphs.Item(4).Range.Start
$rangeEnd = $d.Paragraphs.Item(12).Range.End
$addrRange = $d.Range($rangeStart, $rangeEnd)
$addrRange.Delete()

# ---------------------------------------------------------------------------
# 3) "Bern, " + "{{date}}" runs merge into a single run's text.
# ---------------------------------------------------------------------------
$dateRange = $d.Paragraphs.Item(6).Range
$dateRange.Find.Execute("Bern, {{date}}", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Bern, {{date}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Sehr geehrte {{identity.localized_salutation}} {{identity.localized_
#    title}} {{identity.last_name}}" -> "{{identity.greeting_salutation_and_
#    name}}"
# ---------------------------------------------------------------------------
$greetPara = $d.Paragraphs.Item(12)
$greetRange = $greetPara.Range
$greetRange.Find.Execute( `
    "Sehr geehrte {{identity.localized_salutation}} {{identity.localized_title}} {{identity.last_name}}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{{identity.greeting_salutation_and_name}}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Style tweaks: Normal style suppresses auto-hyphenation, and the
#    "FootnoteAnchor" character style is renamed to "Footnote Reference".
# ---------------------------------------------------------------------------
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = $false

$footnoteAnchor = $d.Styles.Item("FootnoteAnchor")
$footnoteAnchor.NameLocal = "Footnote Reference"

Write-Output "done"
